# Generate Report for Archive
# Two files ("94e2db08-1df9-4b1a-ba22-649f44fbd679.md" and
# "9f05ad46-872c-4edf-baae-8ed791cd5b1c.md") have moved out of the
# "Ready for handoff" state and are now "In Translation". Update the
# Status column on the per-language sheets and the corresponding
# zh-cn / de-de columns on the Overview sheet to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Per-language sheets: column C = "Status", rows 3 and 4 are the two
# files that are now in translation.
foreach ($langSheet in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($langSheet)
    $ws.Cells.Item(3, 3).Value = $newStatus
    $ws.Cells.Item(4, 3).Value = $newStatus
}

# Overview sheet: column E = "zh-cn", column F = "de-de", rows 3 and 4
# are the same two files.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(3, 5).Value = $newStatus
$overview.Cells.Item(3, 6).Value = $newStatus
$overview.Cells.Item(4, 5).Value = $newStatus
$overview.Cells.Item(4, 6).Value = $newStatus
